$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style_D2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.910.65"
$ws.Range("D2").Style = $style_D2
$ws.Range("E2").Value = "  -3.02%  "
$style_D3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.848.22"
$ws.Range("D3").Style = $style_D3
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("E4").Value = "  +0.16%  "
$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.81"
$ws.Range("D5").Style = $style_D5
$ws.Range("E5").Value = "  +0.16%  "
$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.29"
$ws.Range("D6").Style = $style_D6
$ws.Range("E6").Value = "  +5.19%  "
$style_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.666"
$ws.Range("D7").Style = $style_D7
$ws.Range("E7").Value = "  -2.37%  "
$ws.Range("E8").Value = "  +0.26%  "
$style_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.743"
$ws.Range("D9").Style = $style_D9
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("E10").Value = "  +3.88%  "
$style_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.24"
$ws.Range("D11").Style = $style_D11
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E12").Value = "  +0.52%  "
$style_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.26"
$ws.Range("D13").Style = $style_D13
$ws.Range("E13").Value = "  +2.72%  "
$style_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.462.67"
$ws.Range("D14").Style = $style_D14
$ws.Range("E14").Value = "  -3.41%  "
$style_D15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.07"
$ws.Range("D15").Style = $style_D15
$ws.Range("E15").Value = "  +3.65%  "
$style_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.855.48"
$ws.Range("D16").Style = $style_D16
$ws.Range("E16").Value = "  -3.29%  "
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("E18").Value = "  -5.00%  "
$ws.Range("E19").Value = "  -2.30%  "
$style_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.816.16"
$ws.Range("D20").Style = $style_D20
$ws.Range("E20").Value = "  -2.57%  "
$style_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.17"
$ws.Range("D21").Style = $style_D21
$ws.Range("E21").Value = "  +0.17%  "
$style_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.77"
$ws.Range("D22").Style = $style_D22
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("E23").Value = "  -2.08%  "
$ws.Range("E24").Value = "  -4.13%  "
$style_D25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.81"
$ws.Range("D25").Style = $style_D25
$ws.Range("E25").Value = "  -2.94%  "
$style_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.68"
$ws.Range("D26").Style = $style_D26
$ws.Range("E26").Value = "  +4.72%  "
$style_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.98"
$ws.Range("D27").Style = $style_D27
$ws.Range("E27").Value = "  -8.42%  "
$style_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.96"
$ws.Range("D28").Style = $style_D28
$ws.Range("E28").Value = "  +0.17%  "
$style_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.49"
$ws.Range("D29").Style = $style_D29
$ws.Range("E29").Value = "  +0.34%  "
$style_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.87"
$ws.Range("D30").Style = $style_D30
$ws.Range("E30").Value = "  -3.92%  "
$style_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.28"
$ws.Range("D31").Style = $style_D31
$ws.Range("E31").Value = "  +6.12%  "
$style_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.47"
$ws.Range("D32").Style = $style_D32
$ws.Range("E32").Value = "  -1.81%  "
$style_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "47.66"
$ws.Range("D33").Style = $style_D33
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("E34").Value = "  -4.21%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$style_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "69.01"
$ws.Range("D35").Style = $style_D35
$ws.Range("E35").Value = "  -2.70%  "
$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "643.08"
$ws.Range("D36").Style = $style_D36
$ws.Range("E36").Value = "  -3.48%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0979"
$ws.Range("E37").Value = "  +8.36%  "
$style_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.436"
$ws.Range("D38").Style = $style_D38
$ws.Range("E38").Value = "  -0.40%  "
$style_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.148"
$ws.Range("D39").Style = $style_D39
$ws.Range("E39").Value = "  +1.66%  "
$style_D40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = $style_D40
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = $style_D41
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$style_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.97"
$ws.Range("D42").Style = $style_D42
$ws.Range("E42").Value = "  +14.02%  "
$style_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.24"
$ws.Range("D43").Style = $style_D43
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("E44").Value = "  +23.15%  "
$style_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0469"
$ws.Range("D45").Style = $style_D45
$ws.Range("E45").Value = "  -4.09%  "
$ws.Range("E46").Value = "  -5.57%  "
$ws.Range("E47").Value = "  -4.19%  "
$style_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.912.02"
$ws.Range("D48").Style = $style_D48
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$style_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.30"
$ws.Range("D49").Style = $style_D49
$ws.Range("E49").Value = "  -3.97%  "
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$style_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.80"
$ws.Range("D50").Style = $style_D50
$ws.Range("E50").Value = "  -16.64%  "
$style_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000279"
$ws.Range("D51").Style = $style_D51
$ws.Range("E51").Value = "  +3.57%  "
